$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix name order for Niko Beerenwinkel: was stored as last="Niko"? No -
# previously A31="Beerenwinkel" (first col) and B31="Niko" (last col).
# Correct order is first="Niko", last="Beerenwinkel".
$ws.Range("A31").Value = "Niko"
$ws.Range("B31").Value = "Beerenwinkel"

# Katherine Lee: institution Murdoch -> MCRI
$ws.Range("C32").Value = "MCRI"

# Paul Gustafson: country USA -> CAN
$ws.Range("D33").Value = "CAN"

# Kelly Van Lancker: country NL -> BE
$ws.Range("D37").Value = "BE"

# Samuel Muller: institution U Macquarie -> Macquarie U
$ws.Range("C40").Value = "Macquarie U"
